# Added sprites for existing BT:
# Insert a new paragraph right after the "scene management and checking"
# hyperlink paragraph, containing a hyperlink to the pngmaker.ai sun
# prompt followed by a short note, matching the diary's existing
# "<link> - <note>" pattern.

$d = $word.ActiveDocument

# Locate the paragraph whose text contains the scene-management link's
# trailing note so the insertion point is robust to absolute index drift.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*scene management and checking*") {
        $anchor = $para
        break
    }
}

# Create a brand-new empty paragraph immediately after the anchor.
$newPara = $anchor.Range.InsertParagraphAfter()

# Re-fetch the paragraph that now occupies the anchor's old "next" slot -
# that is our freshly inserted (still empty) paragraph.
$target = $anchor.Next().Range
$startPos = $target.Start

$url = "https://pngmaker.ai/app?prompt=sun"
$suffix = " " + [char]0x2013 + " ai png maker used for tempory images of BT"

# Insert the full text first ...
$target.InsertAfter($url + $suffix)

# ... then turn just the URL portion into a real hyperlink, leaving the
# trailing note as plain text, mirroring the rest of the diary entries.
$urlRange = $d.Range($startPos, $startPos + $url.Length)
$d.Hyperlinks.Add($urlRange, $url)
